$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 (the empty row between "name" header and "Belgium test" row),
# which shifts rows 3-6 up to rows 2-5.
$ws.Rows("2").Delete()

# Update the active selection to D6 as in the target file
$ws.Range("D6").Select()
